$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the ISBN header text in F1 (was an empty shared-string placeholder)
$ws.Range("F1").Value = "ISBN"

# Remove the 3 sample data rows (rows 2-4), leaving just the header row
$ws.Rows("2:4").Delete()
